$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.016.40'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.910.07'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("D4").Value = '''0.9999'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '''0.7946'
$ws.Range("E5").Value = '  +5.78%  '
$ws.Range("D6").Value = '''242.09'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").Value = '''0.3169'
$ws.Range("E8").Value = '  +3.71%  '
$ws.Range("D9").Value = '''26.43'
$ws.Range("E9").Value = '  +5.06%  '
$ws.Range("D10").Value = '''0.06944'
$ws.Range("E10").Value = '  +1.59%  '
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = '1.908.40'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '''0.7441'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").Value = '''5.199'
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("D15").Value = '''93.15'
$ws.Range("E15").Value = '  +2.16%  '
$ws.Range("D16").Value = '30.011.34'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").Value = '''14.00'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '''5.883'
$ws.Range("E18").Value = '  -4.37%  '
$ws.Range("D19").Value = '''246.86'
$ws.Range("E19").Value = '  +4.96%  '
$ws.Range("D20").Value = '''0.000007767'
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = '2.151.96'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").Value = '''6.850'
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("D25").Value = '''168.33'
$ws.Range("E25").Value = '  +1.58%  '
$ws.Range("D26").Value = '''9.244'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '''0.1393'
$ws.Range("E27").Value = '  +7.89%  '
$ws.Range("D28").Value = '''18.94'
$ws.Range("E28").Value = '  +1.44%  '
$ws.Range("D29").Value = '''2.035'
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("D30").Value = '''1.365'
$ws.Range("E30").Value = '  +2.00%  '
$ws.Range("D31").Value = '''1.514'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D32").Value = '''4.319'
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.05585'
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '''4.093'
$ws.Range("E34").Value = '  +1.64%  '
$ws.Range("D35").Value = '''1.262'
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("D36").Value = '''0.7334'
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").Value = '''2.718'
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("E38").Value = '  -0.49%  '
$ws.Range("E39").Value = '  +0.99%  '
$ws.Range("D40").Value = '''6.129'
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = '''72.53'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("D44").Value = '''0.8339'
$ws.Range("E44").Value = '  +0.83%  '
$ws.Range("D45").Value = '''1.881'
$ws.Range("E45").Value = '  -2.45%  '
$ws.Range("D46").Value = '''100.65'
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("D47").Value = '''7.557'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").Value = '''985.77'
$ws.Range("E48").Value = '  +7.79%  '
$ws.Range("D49").Value = '2.061.35'
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").Value = '''36.34'
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").Value = '''2.812'
$ws.Range("E51").Value = '  +7.18%  '
